# Fruta / hortaliza, semanal
# Updates columns D (Fecha), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg)
# for rows 2-21 with the refreshed weekly price-report figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @(D, I, J, K, L, M, P)
$data = @{
  2  = @(44544, "Primera", 1000,  600,  650,  625,  625)
  3  = @(44201, "Segunda",  500,  800,  900,  850,  850)
  4  = @(44874, "Tercera", 1200,  450,  500,  475,  475)
  5  = @(44224, "Segunda",  800,  850,  900,  875,  875)
  6  = @(44573, "Tercera",  800,  600,  650,  625,  625)
  7  = @(44278, "Segunda",  700,  600,  700,  650,  650)
  8  = @(44278, "Tercera",  400,  500,  600,  550,  550)
  9  = @(44229, "Segunda",  760,  550,  600,  575,  575)
  10 = @(44799, "Primera",  800, 1000, 1200, 1100, 1100)
  11 = @(44210, "Segunda",  900,  600,  700,  650,  650)
  12 = @(44174, "Segunda",  800,  450,  500,  475,  475)
  13 = @(44174, "Tercera", 1200,  250,  350,  300,  300)
  14 = @(44474, "Segunda",  200,  600,  700,  650,  650)
  15 = @(44253, "Segunda", 1000,  800,  900,  850,  850)
  16 = @(44253, "Tercera",  800,  600,  700,  650,  650)
  17 = @(44267, "Tercera",  400,  500,  600,  550,  550)
  18 = @(44658, "Segunda", 1000,  600,  650,  625,  625)
  19 = @(44245, "Primera",  800,  850,  900,  875,  875)
  20 = @(44245, "Segunda", 1000,  750,  800,  775,  775)
  21 = @(44935, "Segunda", 1000,  400,  500,  460,  460)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 9).Value  = $vals[1]   # I - Calidad
    $ws.Cells.Item($row, 10).Value = $vals[2]   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[3]   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[4]   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[5]   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[6]   # P - Precio $/Kg
}
